$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 59.68504933333333
$ws.Range("H2").Value = 179.055148
$ws.Range("I2").Value = 0.207862575863973
$ws.Range("J2").Value = 0.2078625758639731
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 25.81052333333333
$ws.Range("N2").Value = 77.43157
$ws.Range("O2").Value = 0.804688548562198
$ws.Range("P2").Value = 0.804688548562198
$ws.Range("Q2").Value = 1540.502358469151
$ws.Range("R2").Value = 13864.52122622236
$ws.Range("S2").Value = 0.1672646344723802
$ws.Range("T2").Value = 0.1672646344723802

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 59.68504933333333
$ws.Range("H3").Value = 179.055148
$ws.Range("I3").Value = 0.207862575863973
$ws.Range("J3").Value = 0.2078625758639731
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.1747726666666667
$ws.Range("N3").Value = 0.5243180000000001
$ws.Range("O3").Value = 0.005448845870037694
$ws.Range("P3").Value = 0.005448845870037694
$ws.Range("Q3").Value = 10.43131523211822
$ws.Range("R3").Value = 93.881837089064
$ws.Range("S3").Value = 0.001132611138031806
$ws.Range("T3").Value = 0.001132611138031806

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 59.68504933333333
$ws.Range("H4").Value = 179.055148
$ws.Range("I4").Value = 0.207862575863973
$ws.Range("J4").Value = 0.2078625758639731
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.089875666666667
$ws.Range("N4").Value = 18.269627
$ws.Range("O4").Value = 0.1898626055677645
$ws.Range("P4").Value = 0.1898626055677645
$ws.Range("Q4").Value = 363.4745295988662
$ws.Range("R4").Value = 3271.270766389796
$ws.Range("S4").Value = 0.03946533025356103
$ws.Range("T4").Value = 0.03946533025356103

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 186.4134573333333
$ws.Range("H5").Value = 559.240372
$ws.Range("I5").Value = 0.6492141976897894
$ws.Range("J5").Value = 0.6492141976897894
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.81052333333333
$ws.Range("N5").Value = 77.43157
$ws.Range("O5").Value = 0.804688548562198
$ws.Range("P5").Value = 0.804688548562198
$ws.Range("Q5").Value = 4811.428890149337
$ws.Range("R5").Value = 43302.86001134403
$ws.Range("S5").Value = 0.5224152304449685
$ws.Range("T5").Value = 0.5224152304449685

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 186.4134573333333
$ws.Range("H6").Value = 559.240372
$ws.Range("I6").Value = 0.6492141976897894
$ws.Range("J6").Value = 0.6492141976897894
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.1747726666666667
$ws.Range("N6").Value = 0.5243180000000001
$ws.Range("O6").Value = 0.005448845870037694
$ws.Range("P6").Value = 0.005448845870037694
$ws.Range("Q6").Value = 32.57997704069956
$ws.Range("R6").Value = 293.219793366296
$ws.Range("S6").Value = 0.003537468099851844
$ws.Range("T6").Value = 0.003537468099851844

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 186.4134573333333
$ws.Range("H7").Value = 559.240372
$ws.Range("I7").Value = 0.6492141976897894
$ws.Range("J7").Value = 0.6492141976897894
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.089875666666667
$ws.Range("N7").Value = 18.269627
$ws.Range("O7").Value = 0.1898626055677645
$ws.Range("P7").Value = 0.1898626055677645
$ws.Range("Q7").Value = 1135.234777753471
$ws.Range("R7").Value = 10217.11299978124
$ws.Range("S7").Value = 0.1232614991449691
$ws.Range("T7").Value = 0.1232614991449691

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 41.03855533333333
$ws.Range("H8").Value = 123.115666
$ws.Range("I8").Value = 0.1429232264462375
$ws.Range("J8").Value = 0.1429232264462375
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.81052333333333
$ws.Range("N8").Value = 77.43157
$ws.Range("O8").Value = 0.804688548562198
$ws.Range("P8").Value = 0.804688548562198
$ws.Range("Q8").Value = 1059.226589997291
$ws.Range("R8").Value = 9533.03930997562
$ws.Range("S8").Value = 0.1150086836448492
$ws.Range("T8").Value = 0.1150086836448492

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 41.03855533333333
$ws.Range("H9").Value = 123.115666
$ws.Range("I9").Value = 0.1429232264462375
$ws.Range("J9").Value = 0.1429232264462375
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.1747726666666667
$ws.Range("N9").Value = 0.5243180000000001
$ws.Range("O9").Value = 0.005448845870037694
$ws.Range("P9").Value = 0.005448845870037694
$ws.Range("Q9").Value = 7.172417751754224
$ws.Range("R9").Value = 64.551759765788
$ws.Range("S9").Value = 0.0007787666321540432
$ws.Range("T9").Value = 0.0007787666321540432

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 41.03855533333333
$ws.Range("H10").Value = 123.115666
$ws.Range("I10").Value = 0.1429232264462375
$ws.Range("J10").Value = 0.1429232264462375
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 6.089875666666667
$ws.Range("N10").Value = 18.269627
$ws.Range("O10").Value = 0.1898626055677645
$ws.Range("P10").Value = 0.1898626055677645
$ws.Range("Q10").Value = 249.9196995196203
$ws.Range("R10").Value = 2249.277295676582
$ws.Range("S10").Value = 0.02713577616923427
$ws.Range("T10").Value = 0.02713577616923427

Write-Host "applied updates"
